$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> index helper (B=2, C=3, D=4, E=5)
$colB = 2
$colC = 3
$colD = 4
$colE = 5

function Set-Cell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    # Force text storage so numeric-looking strings (trailing zeros,
    # percent signs, etc.) are preserved verbatim, matching the source
    # workbook which stores these columns as literal text.
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2 (BNB)
Set-Cell 2 $colD "245.22"
Set-Cell 2 $colE "-0.72%"

# Row 3 (OKB)
Set-Cell 3 $colD "27.19"
Set-Cell 3 $colE "2.81%"

# Row 4 (HuobiToken)
Set-Cell 4 $colD "5.109"
Set-Cell 4 $colE "0.41%"

# Row 5 (Cronos)
Set-Cell 5 $colD "0.05705"
Set-Cell 5 $colE "1.88%"

# Row 6 (KuCoinToken)
Set-Cell 6 $colD "6.501"

# Row 7 (MXToken)
Set-Cell 7 $colD "0.8193"
Set-Cell 7 $colE "0.73%"

# Row 8 (FTXToken)
Set-Cell 8 $colD "0.8574"
Set-Cell 8 $colE "1.55%"

# Row 9 -> MandalaExchangeToken
Set-Cell 9 $colB "MandalaExchangeToken"
Set-Cell 9 $colC "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-Cell 9 $colD "0.06946"
Set-Cell 9 $colE "-0.50%"

# Row 10 -> BitrueCoin
Set-Cell 10 $colB "BitrueCoin"
Set-Cell 10 $colC "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-Cell 10 $colD "0.02853"
Set-Cell 10 $colE "0.22%"

# Row 11 -> BitMartToken
Set-Cell 11 $colB "BitMartToken"
Set-Cell 11 $colC "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-Cell 11 $colD "0.09399"
Set-Cell 11 $colE "0.01%"

# Row 12 -> BitForexToken
Set-Cell 12 $colB "BitForexToken"
Set-Cell 12 $colC "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-Cell 12 $colD "0.001517"
Set-Cell 12 $colE "0.55%"

# Row 13 -> CoinExToken
Set-Cell 13 $colB "CoinExToken"
Set-Cell 13 $colC "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-Cell 13 $colD "0.04023"
Set-Cell 13 $colE "-13.55%"

# Row 14 -> One
Set-Cell 14 $colB "One"
Set-Cell 14 $colC "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-Cell 14 $colD "0.0005978"
Set-Cell 14 $colE "0.34%"

# Row 15 -> TigerCash
Set-Cell 15 $colB "TigerCash"
Set-Cell 15 $colC "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-Cell 15 $colD "0.006147"
Set-Cell 15 $colE "-1.37%"

# Row 16 -> LEO
Set-Cell 16 $colB "LEO"
Set-Cell 16 $colC "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-Cell 16 $colD "3.512"
Set-Cell 16 $colE "-2.77%"

# Row 17 -> GateToken
Set-Cell 17 $colB "GateToken"
Set-Cell 17 $colC "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-Cell 17 $colD "3.007"
Set-Cell 17 $colE "-0.24%"

# Row 18 -> BTSEToken
Set-Cell 18 $colB "BTSEToken"
Set-Cell 18 $colC "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-Cell 18 $colD "2.230"
Set-Cell 18 $colE "8.50%"

# Row 19 -> BitpandaEcosystemToken
Set-Cell 19 $colB "BitpandaEcosystemToken"
Set-Cell 19 $colC "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-Cell 19 $colD "0.3165"
Set-Cell 19 $colE "-0.37%"

# Row 20 -> WazirX
Set-Cell 20 $colB "WazirX"
Set-Cell 20 $colC "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-Cell 20 $colD "0.1329"
Set-Cell 20 $colE "-0.49%"

# Row 21 (LiechtensteinCryptoassetsExchange)
Set-Cell 21 $colE "1.45%"

# Row 22 (ProBitToken)
Set-Cell 22 $colE "0.46%"

# Row 23 (MCDex)
Set-Cell 23 $colD "3.567"
Set-Cell 23 $colE "-4.56%"

# Row 24 (ZBToken)
Set-Cell 24 $colE "1.78%"

# Row 25 (BitKan)
Set-Cell 25 $colD "0.001218"
Set-Cell 25 $colE "-2.28%"

# Row 26 (HotbitToken)
Set-Cell 26 $colD "0.004474"

# Row 27 (NitroEx)
Set-Cell 27 $colD "0.00009896"
Set-Cell 27 $colE "3.11%"

# Row 28 (UpBots)
Set-Cell 28 $colE "3.68%"

# Row 40 (IDEX)
Set-Cell 40 $colD "0.03727"
Set-Cell 40 $colE "1.52%"

# Row 41 (KickToken)
Set-Cell 41 $colD "0.005938"
Set-Cell 41 $colE "-3.90%"

# Row 42 (BKEXToken)
Set-Cell 42 $colD "0.1058"
Set-Cell 42 $colE "0.03%"

# Row 43 (CEJI)
Set-Cell 43 $colD "0.002469"
Set-Cell 43 $colE "-1.22%"

# Row 44 (LocalTraders)
Set-Cell 44 $colD "0.009547"
Set-Cell 44 $colE "6.82%"

# Row 45 (CoinLion)
Set-Cell 45 $colD "0.00005137"
Set-Cell 45 $colE "-4.19%"

# Row 46 (Kangarootoken)
Set-Cell 46 $colE "-0.02%"

# Row 47 (CoinbaseStockToken)
Set-Cell 47 $colE "-8.18%"

# Row 48 (BOLO)
Set-Cell 48 $colD "0.002517"
Set-Cell 48 $colE "-3.74%"

# Row 49 (CryptobidCoin)
Set-Cell 49 $colE "-0.02%"

# Row 50 (SpecialPowerGold)
Set-Cell 50 $colE "-0.02%"
